$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "stimFile" column is being inserted right after the posFile
# column (the old B/C columns - cs_plus/cs_minus file names - shift one
# column to the right, becoming C/D). Inserting a real column (instead
# of just overwriting values in place) lets the existing column-width
# cache for the untouched columns carry over unchanged.
$ws.Columns.Item(2).Insert()

# Clear all cell content so the previously shifted row 3 data (now
# orphaned) is wiped out along with everything else, then rewrite the
# table from scratch.
$ws.Cells.Clear()

# Re-establish the values that are reused unchanged from the original
# sheet first so that they keep their original shared-string order,
# then add the genuinely new values afterwards (matching the order in
# which the new strings were introduced in the target workbook).
$ws.Range("C2").Value = "stimuli/social/016_y_m_n_b.jpg"
$ws.Range("D2").Value = "stimuli/social/031_y_m_n_a.jpg"
$ws.Range("F2").Value = "stimuli/non-social/016_y_m_n_b_scrambled.jpg"
$ws.Range("E2").Value = "stimuli/non-social/031_y_m_n_a_scrambled.jpg"
$ws.Range("A1").Value = "posFile"
$ws.Range("A2").Value = "positions.xlsx"

# New values introduced by this edit
$ws.Range("C1").Value = "cs_plus_s"
$ws.Range("D1").Value = "cs_minus_s"
$ws.Range("E1").Value = "cs_plus_ns"
$ws.Range("F1").Value = "cs_minus_ns"
$ws.Range("B1").Value = "stimFile"
$ws.Range("B2").Value = "stimuli.xlsx"

# New column B takes on the same (roughly 16.5 character) width as
# column A.
$ws.Columns.Item(2).ColumnWidth = 15.67

# Selection
$ws.Range("B1:B2").Select()
